$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F2").Value = 1641
$ws1.Range("F3").Value = 861
$ws1.Range("F4").Value = 267
$ws1.Range("F5").Value = 82
$ws1.Range("F6").Value = 1180
$ws1.Range("F7").Value = 792
$ws1.Range("F8").Value = 825
$ws1.Range("F9").Value = 1516
$ws1.Range("F10").Value = 302
$ws1.Range("F11").Value = 1055
$ws1.Range("F13").Value = 74
$ws1.Range("F14").Value = 201
$ws1.Range("F15").Value = 58
$ws1.Range("F16").Value = 509
$ws1.Range("F17").Value = 60
$ws1.Range("F18").Value = 40
$ws1.Range("F21").Value = 302
$ws1.Range("F22").Value = 577
$ws1.Range("F23").Value = 583
$ws1.Range("F24").Value = 50
$ws1.Range("F25").Value = 8
$ws1.Range("F26").Value = 777
$ws1.Range("F27").Value = 259
$ws1.Range("F28").Value = 194

$ws2.Range("F3").Value = 1030
$ws2.Range("F5").Value = 280
$ws2.Range("F9").Value = 596
$ws2.Range("F10").Value = 88

$ws4.Range("F3").Value = 1641
$ws4.Range("F5").Value = 861
$ws4.Range("F6").Value = 267
$ws4.Range("F7").Value = 1030
$ws4.Range("F8").Value = 82
$ws4.Range("F9").Value = 1180
$ws4.Range("F10").Value = 792
$ws4.Range("F11").Value = 825
$ws4.Range("F12").Value = 1517
$ws4.Range("F13").Value = 302
$ws4.Range("F14").Value = 1055
$ws4.Range("F16").Value = 74
$ws4.Range("F17").Value = 201
$ws4.Range("F18").Value = 58
$ws4.Range("F19").Value = 509
$ws4.Range("F20").Value = 60
$ws4.Range("F21").Value = 40
$ws4.Range("F24").Value = 280
$ws4.Range("F26").Value = 302
$ws4.Range("F30").Value = 577
$ws4.Range("F31").Value = 583
$ws4.Range("F32").Value = 50
$ws4.Range("F33").Value = 8
$ws4.Range("F34").Value = 777
$ws4.Range("F35").Value = 259
$ws4.Range("F37").Value = 194
$ws4.Range("F38").Value = 596
$ws4.Range("F39").Value = 88
$ws4.Range("F40").Value = 88
